# Small improvements to plotting vignette
# - Fix typo "Aciclovr2" -> "Aciclovir2" on the plotGrids sheet
# - Add a new exportConfiguration row for Aciclovir.png
# - Update active sheet / selections to match the author's final view

$wb = $excel.ActiveWorkbook

# --- plotGrids: fix typo in A3 ("Aciclovr2" -> "Aciclovir2") ---
$wsPlotGrids = $wb.Worksheets.Item("plotGrids")
$wsPlotGrids.Range("A3").Value = "Aciclovir2"
$wsPlotGrids.Range("A2").Select()

# --- exportConfiguration: add new row with Aciclovir / Aciclovir.png ---
$wsExportConfig = $wb.Worksheets.Item("exportConfiguration")
$wsExportConfig.Range("A2").Value = "Aciclovir"
$wsExportConfig.Range("B2").Value = "Aciclovir.png"

# Make exportConfiguration the active sheet with B2 selected
$wsExportConfig.Activate()
$wsExportConfig.Range("B2").Select()
